# Updated symbol list: refresh Price (D) and Volume(1h) (E) columns with latest
# cryptocurrency market data, matching the Jan 24 2023 GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These columns are stored as plain text (e.g. "303.32", "-1.08%") rather than
# numbers/percentages, so we force a Text number format while writing the new
# values and then clear the formatting again to avoid altering cell appearance.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$changes = @(
    @{Cell="D2"; Value="303.32"},
    @{Cell="E2"; Value="-1.08%"},
    @{Cell="D3"; Value="35.27"},
    @{Cell="E3"; Value="-2.34%"},
    @{Cell="D4"; Value="5.011"},
    @{Cell="E4"; Value="-1.24%"},
    @{Cell="D5"; Value="0.07839"},
    @{Cell="E5"; Value="-1.46%"},
    @{Cell="D6"; Value="1.870"},
    @{Cell="E6"; Value="-14.02%"},
    @{Cell="D7"; Value="4.089"},
    @{Cell="E7"; Value="-1.68%"},
    @{Cell="D8"; Value="7.813"},
    @{Cell="E8"; Value="-2.68%"},
    @{Cell="E9"; Value="8.10%"},
    @{Cell="D10"; Value="0.9211"},
    @{Cell="E10"; Value="-1.03%"},
    @{Cell="D11"; Value="0.1064"},
    @{Cell="E11"; Value="8.45%"},
    @{Cell="D12"; Value="0.1870"},
    @{Cell="E12"; Value="-0.49%"},
    @{Cell="D13"; Value="0.09408"},
    @{Cell="E13"; Value="4.06%"},
    @{Cell="D14"; Value="0.03619"},
    @{Cell="E14"; Value="0.05%"},
    @{Cell="D15"; Value="0.09949"},
    @{Cell="E15"; Value="0.30%"},
    @{Cell="D16"; Value="0.001408"},
    @{Cell="E16"; Value="-1.93%"},
    @{Cell="D17"; Value="0.005732"},
    @{Cell="E17"; Value="0.55%"},
    @{Cell="D18"; Value="3.462"},
    @{Cell="E18"; Value="-0.48%"},
    @{Cell="D19"; Value="0.3431"},
    @{Cell="E19"; Value="1.79%"},
    @{Cell="D20"; Value="0.1295"},
    @{Cell="E20"; Value="-4.60%"},
    @{Cell="E21"; Value="1.14%"},
    @{Cell="D22"; Value="0.2200"},
    @{Cell="E22"; Value="0.31%"},
    @{Cell="D23"; Value="0.04551"},
    @{Cell="E23"; Value="-0.67%"},
    @{Cell="D24"; Value="0.001230"},
    @{Cell="E24"; Value="-0.72%"},
    @{Cell="D25"; Value="0.004657"},
    @{Cell="E25"; Value="-2.26%"},
    @{Cell="D26"; Value="0.0001257"},
    @{Cell="E26"; Value="-3.41%"},
    @{Cell="D27"; Value="0.0004472"},
    @{Cell="E27"; Value="-5.62%"},
    @{Cell="D39"; Value="0.01892"},
    @{Cell="E39"; Value="-3.00%"},
    @{Cell="D40"; Value="0.04727"},
    @{Cell="E40"; Value="-3.81%"},
    @{Cell="D41"; Value="0.007645"},
    @{Cell="E41"; Value="-1.75%"},
    @{Cell="E42"; Value="28.28%"},
    @{Cell="D43"; Value="0.1337"},
    @{Cell="E44"; Value="1.31%"},
    @{Cell="E45"; Value="-1.59%"},
    @{Cell="D46"; Value="0.00006279"},
    @{Cell="E46"; Value="1.03%"},
    @{Cell="D47"; Value="0.00000000754"},
    @{Cell="E47"; Value="0.46%"},
    @{Cell="D48"; Value="66.41"},
    @{Cell="E48"; Value="27.74%"},
    @{Cell="D49"; Value="0.001307"},
    @{Cell="E49"; Value="-27.45%"},
    @{Cell="D50"; Value="0.00002112"},
    @{Cell="E50"; Value="0.46%"},
    @{Cell="D51"; Value="0.0002011"},
    @{Cell="E51"; Value="0.46%"}
)

foreach ($chg in $changes) {
    $ws.Range($chg.Cell).Value = $chg.Value
}

$dataRange.ClearFormats()
